$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 09:22"

# --- Row 31: Polonia -- update Casos activos / Recuperados ---
$ws.Range("D31").Value = 1297
$ws.Range("E31").Value = 7916

# --- Rows 71-75: Oman inserted before Nueva Zelanda, Armenia/Lituania swap+update ---
# Row 71 becomes Oman (new data)
$ws.Range("A71").Value = "Oman"
$ws.Range("B71").Value = 1508
$ws.Range("C71").Value = 98
$ws.Range("D71").Value = 238
$ws.Range("E71").Value = 1262
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 8

# Row 72 becomes Nueva Zelanda (old row71 data, unchanged values)
$ws.Range("A72").Value = "Nueva Zelanda"
$ws.Range("B72").Value = 1445
$ws.Range("C72").Value = 5
$ws.Range("D72").Value = 1006
$ws.Range("E72").Value = 426
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 13

# Row 73 becomes Azerbaiyan (old row72 data, unchanged values)
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 1436
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 791
$ws.Range("E73").Value = 626
$ws.Range("F73").Value = 16
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 19

# Row 74 becomes Armenia (new data)
$ws.Range("A74").Value = "Armenia"
$ws.Range("B74").Value = 1401
$ws.Range("C74").Value = 62
$ws.Range("D74").Value = 609
$ws.Range("E74").Value = 768
$ws.Range("F74").Value = 30
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 24

# Row 75 becomes Lituania (new data)
$ws.Range("A75").Value = "Lituania"
$ws.Range("B75").Value = 1350
$ws.Range("C75").Value = 24
$ws.Range("D75").Value = 298
$ws.Range("E75").Value = 1015
$ws.Range("F75").Value = 17
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 37

# --- Row 114: Montenegro ---
$ws.Range("B114").Value = 313
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 90
$ws.Range("E114").Value = 218

# --- Row 115: Sri Lanka ---
$ws.Range("B115").Value = 309
$ws.Range("C115").Value = 5
$ws.Range("E115").Value = 204

# --- Row 121: Vietnam ---
$ws.Range("D121").Value = 216
$ws.Range("E121").Value = 52

# --- Row 126: Jamaica ---
$ws.Range("E126").Value = 190
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 6

# --- Row 127: Paraguay ---
$ws.Range("D127").Value = 53
$ws.Range("E127").Value = 147

# --- Row 144: Guayana Francesa ---
$ws.Range("D144").Value = 76
$ws.Range("E144").Value = 20
$ws.Range("F144").Value = 1
